# Daily attendance processing - 2025-10-06 10:40:17
# Applies the day's attendance-recording updates to the session analysis sheet:
#  - newly recorded sessions (rows 15, 86, 140) flip from Pending -> Recorded,
#    picking up the green "Recorded" row style, a "Recorded By" list, a
#    Students count, and a Status of "Recorded"
#  - "Recorded By" lists get re-ordered (same contributors, re-sorted) and in
#    a couple of cases gain an extra contributor
#  - Students counts (H column) tick up for a few already-recorded sessions
#  - the derived Class/Group statistics (K:S columns) are recomputed to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a percentage-looking value ("41.2%") as literal TEXT rather
# than letting Excel's smart-entry reinterpret it as a numeric percentage
# (which would change both the stored value and the cell's number format /
# style). We briefly force a text format, assign the value, then restore the
# original look by copying formatting from an unrelated, never-edited cell
# that already carries the sheet's standard "stat cell" style.
function Set-PercentText($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("K9").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# Rows 15 / 86 / 140: sessions that just got recorded.
# Copy the formatting (fill/style) from an already-"Recorded" row (row 7,
# which uses the green "Recorded" style) onto each newly recorded row, then
# fill in the Recorded By / Students / Status cells and the recomputed
# statistics for row 15.
# ---------------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A15:I15").PasteSpecial(-4122)
$ws.Range("A86:I86").PasteSpecial(-4122)
$ws.Range("A140:I140").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("G15").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("H15").Value = "5/217"
$ws.Range("I15").Value = "Recorded"

$ws.Range("G86").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("H86").Value = "17/154"
$ws.Range("I86").Value = "Recorded"

$ws.Range("G140").Value = "Sara_nabil@med.asu.edu.eg"
$ws.Range("H140").Value = "12/250"
$ws.Range("I140").Value = "Recorded"

# ---------------------------------------------------------------------------
# "Recorded By" list re-orderings (same people, new order) and additions.
# ---------------------------------------------------------------------------
$ws.Range("G7").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G8").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G17").Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G24").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G25").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G32").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G34").Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G35").Value = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G37").Value = "System, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G45").Value = "backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G51").Value = "yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G52").Value = "abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G54").Value = "System, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G62").Value = "backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G68").Value = "yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G69").Value = "abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G72").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G76").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G85").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G98").Value = "nourhanmohamed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, user@user.com, nourhanhosni@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G102").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G109").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G116").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G119").Value = "shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G126").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G133").Value = "afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G136").Value = "shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G143").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G151").Value = "nourhanmohamed@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System, marian.samir@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Students (H column) counts that increased.
# ---------------------------------------------------------------------------
$ws.Range("H17").Value = "101/217"
$ws.Range("H51").Value = "52/220"
$ws.Range("H98").Value = "94/224"

# ---------------------------------------------------------------------------
# Class Statistics block (K:L, rows 6/8/9/10).
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 59
$ws.Range("L8").Value = 93
Set-PercentText "L9" "38.6%"
Set-PercentText "L10" "47.2%"

# ---------------------------------------------------------------------------
# Group Statistics block (M:S) recomputed per group after today's recording.
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 7
$ws.Range("Q15").Value = 10
Set-PercentText "R15" "41.2%"
Set-PercentText "S15" "54.7%"

Set-PercentText "S17" "44.9%"

$ws.Range("O19").Value = 6
$ws.Range("Q19").Value = 11
Set-PercentText "R19" "35.3%"
Set-PercentText "S19" "51.6%"

Set-PercentText "S20" "48.4%"

$ws.Range("O23").Value = 7
$ws.Range("Q23").Value = 9
Set-PercentText "R23" "41.2%"
Set-PercentText "S23" "15.8%"
